# zero-excel finished, import data
#
# DATA-USER.xlsx:
#  - B2 label changes from the old I_API entity marker to S_USER
#  - the merged title-row banner (was D2:H2, holding the DAO class name)
#    moves to C2:I2 and now shows the friendlier "账号导入" caption
#  - a new "统一标识符" (sigma) column is appended as column I, mirroring
#    column H's per-row formatting, with the same sigma value for both
#    data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the entity marker in B2.
$ws.Range("B2").Value = "S_USER"

# 2. Pull the banner caption out of D2 before we rearrange the merge, then
#    break the old D2:H2 merge apart.
$banner = $ws.Range("D2").Value()
$ws.Range("D2:H2").UnMerge()

# 3. Move the caption into C2 (replacing the DAO class-name text that used
#    to live there) and blank out the old D2 cell.
$ws.Range("C2").Value = $banner
$ws.Range("D2").ClearContents()

# 4. Merge C2:I2 into a single cell first, then stamp every cell of the new
#    banner span with the same fill/font/border look (copied from B2, which
#    already carries the right fill+font+border combination) with
#    left-aligned text. Formatting after the merge keeps every cell in the
#    merged range on one consistent style (merging first and formatting
#    after, rather than the reverse, avoids the merge operation clobbering
#    the formatting we just set).
$ws.Range("C2:I2").Merge()
$ws.Range("B2").Copy()
$ws.Range("C2:I2").PasteSpecial(-4122)
$ws.Range("C2:I2").HorizontalAlignment = -4131

# 5. Add the new "统一标识符" / "sigma" column, reusing column H's
#    per-row formatting (header row, key row, the two data rows).
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = "统一标识符"

$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = "sigma"

$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = "kbm9LQBAsm8BPJQ7AIG9MVDgaF7azrWd"

$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = "kbm9LQBAsm8BPJQ7AIG9MVDgaF7azrWd"

# 6. Resize the columns: the new I column is wide enough for the sigma
#    values, and G narrows now that the sheet has more columns to show.
$ws.Columns.Item(9).ColumnWidth = 45.67
$ws.Columns.Item(7).ColumnWidth = 11

# 7. Leave the selection where the user's cursor ended up after typing the
#    last sigma value.
$ws.Range("I7").Select()
